$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text (row 1)
$ws.Range("A1").Value = "Datos actualizados a 1 de Agosto de 2020 a las 14:15"

# Update Estados Unidos (row 4)
$ws.Range("B4").Value = 4707401
$ws.Range("C4").Value = 1512
$ws.Range("D4").Value = 2328451
$ws.Range("E4").Value = 2222178
$ws.Range("G4").Value = 25
$ws.Range("H4").Value = 156772

# Update India (row 6)
$ws.Range("B6").Value = 1706391
$ws.Range("C6").Value = 9337
$ws.Range("D6").Value = 1100685
$ws.Range("E6").Value = 569044
$ws.Range("G6").Value = 111
$ws.Range("H6").Value = 36662

# Update Catar (row 26)
$ws.Range("B26").Value = 110911
$ws.Range("C26").Value = 216
$ws.Range("D26").Value = 107578
$ws.Range("E26").Value = 3159

# Update Nepal (row 68)
$ws.Range("B68").Value = 20086
$ws.Range("C68").Value = 315
$ws.Range("D68").Value = 14492
$ws.Range("E68").Value = 5538

# Update Madagascar (row 83)
$ws.Range("B83").Value = 11273
$ws.Range("C83").Value = 405
$ws.Range("D83").Value = 8109
$ws.Range("E83").Value = 3057
$ws.Range("G83").Value = 1
$ws.Range("H83").Value = 107

# Update Croacia (row 100)
$ws.Range("B100").Value = 5224
$ws.Range("C100").Value = 86
$ws.Range("D100").Value = 4341
$ws.Range("E100").Value = 738

# Update Islandia (row 131)
$ws.Range("B131").Value = 1893
$ws.Range("C131").Value = 8
$ws.Range("E131").Value = 58

# Update Surinam (row 136)
$ws.Range("B136").Value = 1706
$ws.Range("C136").Value = 56
$ws.Range("D136").Value = 1123
$ws.Range("E136").Value = 557

# Vietnam overtakes Bahamas: row 161 now shows Vietnam's (updated) figures,
# row 162 now shows Bahamas's (previous) figures.
$ws.Range("A161").Value = "Vietnam"
$ws.Range("B161").Value = 586
$ws.Range("C161").Value = 40
$ws.Range("D161").Value = 373
$ws.Range("E161").Value = 210
$ws.Range("G161").Value = 1
$ws.Range("H161").Value = 3

$ws.Range("A162").Value = "Bahamas"
$ws.Range("B162").Value = 574
$ws.Range("C162").Value = 0
$ws.Range("D162").Value = 91
$ws.Range("E162").Value = 469
$ws.Range("G162").Value = 0
$ws.Range("H162").Value = 14
